# The "Code" column (G) in the instruction tables on Sheet1 is removed.
# This deletes the entire column, shifting the "Operandos" column (and
# everything to its right) one position to the left, which is exactly
# what the target OOXML shows (dimension A2:I45 -> A2:H45, the shared
# strings "Code" and "0x60" become unused and drop out on save, etc.).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("G").Delete()

# Match the author's final selection in the saved file.
$ws.Range("E7").Select()
